# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Wed Feb  7 23:28:47 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.205.13'
$ws.Range('E2').Value = '  +2.50%  '
$ws.Range('D3').Value = '2.426.01'
$ws.Range('E3').Value = '  +2.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.94'
$ws.Range('E5').Value = '  +1.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.22'
$ws.Range('E6').Value = '  +4.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.513'
$ws.Range('E7').Value = '  +1.37%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.502'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.37'
$ws.Range('E10').Value = '  +3.46%  '
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('E12').Value = '  +4.11%  '
$ws.Range('E13').Value = '  +2.12%  '
$ws.Range('E14').Value = '  +2.03%  '
$ws.Range('D15').Value = '2.804.46'
$ws.Range('E15').Value = '  +2.02%  '
$ws.Range('D16').Value = '2.430.48'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.837'
$ws.Range('E17').Value = '  +3.72%  '
$ws.Range('D18').Value = '44.168.42'
$ws.Range('E18').Value = '  +2.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.27'
$ws.Range('E19').Value = '  +1.02%  '
$ws.Range('D21').Value = '0.0₃0906'
$ws.Range('E21').Value = '  +2.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.59'
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('E23').Value = '  +5.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.49'
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('E25').Value = '  +2.03%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.31'
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('E28').Value = '  -1.58%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.58'
$ws.Range('E29').Value = '  +4.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.89'
$ws.Range('E30').Value = '  +4.90%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.117'
$ws.Range('E31').Value = '  +11.09%  '
$ws.Range('B32').Value = 'Celestia'
$ws.Range('C32').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.66'
$ws.Range('E32').Value = '  +7.31%  '
$ws.Range('E33').Value = '  +2.47%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  +1.57%  '
$ws.Range('E36').Value = '  +3.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.47'
$ws.Range('E37').Value = '  +4.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '129.42'
$ws.Range('E38').Value = '  +24.36%  '
$ws.Range('E39').Value = '  +3.79%  '
$ws.Range('E40').Value = '  -1.12%  '
$ws.Range('E41').Value = '  +0.78%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.33'
$ws.Range('E42').Value = '  -4.74%  '
$ws.Range('E43').Value = '  +2.40%  '
$ws.Range('D44').Value = '1.952.21'
$ws.Range('E44').Value = '  -0.44%  '
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.88'
$ws.Range('E46').Value = '  +5.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.46'
$ws.Range('E47').Value = '  +3.55%  '
$ws.Range('D49').Value = '2.665.62'
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.53'
$ws.Range('E50').Value = '  +1.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.81'
$ws.Range('E51').Value = '  +2.65%  '

Write-Host "Applied cryptos update"
